$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (avoid Excel auto-converting numeric-looking strings
# like "1.000" / "0.2050" into actual numbers) across the edited block,
# then strip the temporary formatting back off once all values are set so
# cells keep their original (default/general) style, matching the source
# data which stored every value as plain inline text.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = '22.027.48'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").Value = '1.553.25'
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("D6").Value = '290.36'
$ws.Range("E6").Value = '  +1.23%  '

$ws.Range("D7").Value = '0.3925'
$ws.Range("E7").Value = '  +2.59%  '

$ws.Range("D8").Value = '0.3208'
$ws.Range("E8").Value = '  -1.79%  '

$ws.Range("D9").Value = '44.31'
$ws.Range("E9").Value = '  +2.23%  '

$ws.Range("D10").Value = '0.07172'
$ws.Range("E10").Value = '  -1.88%  '

$ws.Range("D11").Value = '1.071'
$ws.Range("E11").Value = '  -4.50%  '

$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.05%  '

$ws.Range("D13").Value = '5.622'
$ws.Range("E13").Value = '  -2.57%  '

$ws.Range("D14").Value = '18.62'
$ws.Range("E14").Value = '  -6.46%  '

$ws.Range("D15").Value = '6.662'
$ws.Range("E15").Value = '  -1.13%  '

$ws.Range("D16").Value = '1.558.60'
$ws.Range("E16").Value = '  -0.80%  '

$ws.Range("D17").Value = '0.00001096'
$ws.Range("E17").Value = '  +1.47%  '

$ws.Range("D18").Value = '0.06563'
$ws.Range("E18").Value = '  -0.76%  '

$ws.Range("D19").Value = '83.23'
$ws.Range("E19").Value = '  -2.51%  '

$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  -0.03%  '

$ws.Range("D21").Value = '6.183'
$ws.Range("E21").Value = '  -2.76%  '

$ws.Range("D22").Value = '15.39'
$ws.Range("E22").Value = '  -3.79%  '

$ws.Range("D23").Value = '11.14'
$ws.Range("E23").Value = '  -4.05%  '

$ws.Range("D24").Value = '22.038.64'
$ws.Range("E24").Value = '  -0.16%  '

$ws.Range("D25").Value = '2.356'
$ws.Range("E25").Value = '  +3.04%  '

$ws.Range("D26").Value = '2.388'
$ws.Range("E26").Value = '  -3.29%  '

$ws.Range("D27").Value = '147.77'
$ws.Range("E27").Value = '  -1.39%  '

$ws.Range("D28").Value = '18.43'
$ws.Range("E28").Value = '  -2.92%  '

$ws.Range("D29").Value = '4.868'
$ws.Range("E29").Value = '  -1.10%  '

$ws.Range("D30").Value = '1.734.45'
$ws.Range("E30").Value = '  -0.48%  '

$ws.Range("D31").Value = '118.03'
$ws.Range("E31").Value = '  -2.00%  '

$ws.Range("D32").Value = '0.9819'
$ws.Range("E32").Value = '  -8.54%  '

$ws.Range("D33").Value = '5.881'
$ws.Range("E33").Value = '  +0.77%  '

$ws.Range("D34").Value = '0.08289'
$ws.Range("E34").Value = '  +1.03%  '

$ws.Range("D35").Value = '9.098'
$ws.Range("E35").Value = '  -0.92%  '

$ws.Range("E36").Value = '  -14.07%  '

$ws.Range("D37").Value = '0.02252'
$ws.Range("E37").Value = '  -2.34%  '

$ws.Range("D38").Value = '5.079'
$ws.Range("E38").Value = '  -2.80%  '

$ws.Range("D39").Value = '0.05996'
$ws.Range("E39").Value = '  -3.66%  '

$ws.Range("D40").Value = '1.203'
$ws.Range("E40").Value = '  -2.33%  '

$ws.Range("D41").Value = '0.2050'
$ws.Range("E41").Value = '  -4.18%  '

$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("D43").Value = '10.70'
$ws.Range("E43").Value = '  -2.20%  '

$ws.Range("D44").Value = '0.5773'
$ws.Range("E44").Value = '  -3.49%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.04'
$ws.Range("E45").Value = '  -3.49%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '3.747'
$ws.Range("E46").Value = '  +0.61%  '

$ws.Range("D47").Value = '0.5542'
$ws.Range("E47").Value = '  -3.91%  '

$ws.Range("D48").Value = '117.01'
$ws.Range("E48").Value = '  -3.67%  '

$ws.Range("D49").Value = '1.871'

$ws.Range("D50").Value = '1.130'
$ws.Range("E50").Value = '  -3.01%  '

$ws.Range("D51").Value = '0.06807'
$ws.Range("E51").Value = '  -2.78%  '

$editRange.ClearFormats()
